$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Workbook window geometry (cosmetic - matches the resized/repositioned
# Excel window recorded in the saved workbookView)
# ---------------------------------------------------------------------------
$excel.Width = 29040
$excel.Height = 15840
$excel.Left = 28680
$excel.Top = -120

# ---------------------------------------------------------------------------
# Sheet1 - update the Year 1 reporting year from 2002 to 2012
# ---------------------------------------------------------------------------
$ws1.Range("C1").Value = "2012"
$ws1.Range("E7").Value = "2012"

# ---------------------------------------------------------------------------
# Sheet1 - updated factor values (Year 1 / "Average Values" column, E) and
# updated Riddership Effect absolute values (column H)
# ---------------------------------------------------------------------------
$ws1.Range("E8").Value = 20275984
$ws1.Range("H8").Value = 4450449.50441

$ws1.Range("E9").Value = 0.53524007
$ws1.Range("H9").Value = 751911.54587

$ws1.Range("E10").Value = 3221377.17
$ws1.Range("H10").Value = 2288619.5879

$ws1.Range("E11").Value = 25.44625052
$ws1.Range("H11").Value = -321272.69507

$ws1.Range("E12").Value = 3.8041
$ws1.Range("H12").Value = -2257292.48741

$ws1.Range("E13").Value = 26799.83
$ws1.Range("H13").Value = -598838.754013

$ws1.Range("E14").Value = 6.77
$ws1.Range("H14").Value = 22313.66406

$ws1.Range("E15").Value = 4.1
$ws1.Range("H15").Value = -27686.239538

# Row 19 "New Reporters" now carries an explicit zero Riddership Effect
$ws1.Range("H19").Value = 0

# Total Modeled / Total Observed ridership (Year 1 column, E)
$ws1.Range("E20").Value = 44561420.18
$ws1.Range("E21").Value = 45966223

# ---------------------------------------------------------------------------
# Sheet1 - the "% Diff" (G) and "% Diff" riddership-share (I) columns drop
# the *100 multiplier now that the cells carry a percentage number format
# ---------------------------------------------------------------------------
for ($r = 8; $r -le 19; $r++) {
    $ws1.Range("G$r").Formula = '=IFERROR((F' + $r + '-E' + $r + ')/E' + $r + ',0)'
    $ws1.Range("I$r").Formula = '=IFERROR(H' + $r + '/$E$21,0)'
}
$ws1.Range("G20").Formula = '=IFERROR((F20-E20)/E20,0)'
$ws1.Range("G21").Formula = '=IFERROR((F21-E21)/E21,0)'
$ws1.Range("I20").Formula = '=G20'
$ws1.Range("I21").Formula = '=G21'

# ---------------------------------------------------------------------------
# Sheet1 - number formats: the Year1/Year2/Riddership-Effect value columns
# (E, F, H) become 2-decimal numbers; the % Diff columns (G, I) become
# 2-decimal percentages
# ---------------------------------------------------------------------------
$ws1.Range("E8:F21").NumberFormat = "#,##0.00"
$ws1.Range("H8:H21").NumberFormat = "#,##0.00"
$ws1.Range("G8:G21").NumberFormat = "0.00%"
$ws1.Range("I8:I21").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# Sheet1 - scroll/selection: the saved view no longer freezes on A7, and the
# active selection moves from K20 to H21
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
$ws1.Range("H21").Select()
